# Applies the Jan 26 2023 "Updated symbol list" crypto price/volume/hour refresh.
# Generated from the canonical OOXML diff: updates Price (D), Volume(1h) (E),
# and Hora (G) columns for rows 2-51 on Sheet1. All of these columns are
# stored as text in the workbook, so each cell's NumberFormat is forced to
# "@" (Text) before assigning the new value, ensuring Excel does not coerce
# the numeric-looking strings (prices, percentages) into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "308.96" },
    @{ Cell = "E2"; Value = "2.32%" },
    @{ Cell = "G2"; Value = "6" },
    @{ Cell = "D3"; Value = "36.19" },
    @{ Cell = "E3"; Value = "3.33%" },
    @{ Cell = "G3"; Value = "6" },
    @{ Cell = "D4"; Value = "5.116" },
    @{ Cell = "E4"; Value = "1.31%" },
    @{ Cell = "G4"; Value = "6" },
    @{ Cell = "D5"; Value = "0.08131" },
    @{ Cell = "E5"; Value = "2.15%" },
    @{ Cell = "G5"; Value = "6" },
    @{ Cell = "D6"; Value = "1.968" },
    @{ Cell = "E6"; Value = "2.79%" },
    @{ Cell = "G6"; Value = "6" },
    @{ Cell = "D7"; Value = "4.184" },
    @{ Cell = "E7"; Value = "3.57%" },
    @{ Cell = "G7"; Value = "6" },
    @{ Cell = "D8"; Value = "7.788" },
    @{ Cell = "E8"; Value = "0.78%" },
    @{ Cell = "G8"; Value = "6" },
    @{ Cell = "D9"; Value = "0.9298" },
    @{ Cell = "E9"; Value = "0.70%" },
    @{ Cell = "G9"; Value = "6" },
    @{ Cell = "D10"; Value = "0.1380" },
    @{ Cell = "E10"; Value = "15.46%" },
    @{ Cell = "G10"; Value = "6" },
    @{ Cell = "D11"; Value = "0.1933" },
    @{ Cell = "E11"; Value = "5.39%" },
    @{ Cell = "G11"; Value = "6" },
    @{ Cell = "D12"; Value = "0.09278" },
    @{ Cell = "E12"; Value = "-0.66%" },
    @{ Cell = "G12"; Value = "6" },
    @{ Cell = "D13"; Value = "0.03407" },
    @{ Cell = "E13"; Value = "-3.65%" },
    @{ Cell = "G13"; Value = "6" },
    @{ Cell = "D14"; Value = "0.09842" },
    @{ Cell = "E14"; Value = "-0.01%" },
    @{ Cell = "G14"; Value = "6" },
    @{ Cell = "D15"; Value = "0.001415" },
    @{ Cell = "E15"; Value = "2.02%" },
    @{ Cell = "G15"; Value = "6" },
    @{ Cell = "D16"; Value = "0.005751" },
    @{ Cell = "E16"; Value = "-1.18%" },
    @{ Cell = "G16"; Value = "6" },
    @{ Cell = "D17"; Value = "3.623" },
    @{ Cell = "G17"; Value = "6" },
    @{ Cell = "D18"; Value = "2.971" },
    @{ Cell = "E18"; Value = "0.58%" },
    @{ Cell = "G18"; Value = "6" },
    @{ Cell = "D19"; Value = "0.3439" },
    @{ Cell = "G19"; Value = "6" },
    @{ Cell = "D20"; Value = "0.1304" },
    @{ Cell = "E20"; Value = "0.95%" },
    @{ Cell = "G20"; Value = "6" },
    @{ Cell = "D21"; Value = "4.887" },
    @{ Cell = "E21"; Value = "-4.06%" },
    @{ Cell = "G21"; Value = "6" },
    @{ Cell = "D22"; Value = "0.2498" },
    @{ Cell = "E22"; Value = "1.29%" },
    @{ Cell = "G22"; Value = "6" },
    @{ Cell = "D23"; Value = "0.04455" },
    @{ Cell = "E23"; Value = "-1.09%" },
    @{ Cell = "G23"; Value = "6" },
    @{ Cell = "E24"; Value = "0.00%" },
    @{ Cell = "G24"; Value = "6" },
    @{ Cell = "D25"; Value = "0.004875" },
    @{ Cell = "E25"; Value = "6.73%" },
    @{ Cell = "G25"; Value = "6" },
    @{ Cell = "E26"; Value = "-0.75%" },
    @{ Cell = "G26"; Value = "6" },
    @{ Cell = "G27"; Value = "6" },
    @{ Cell = "G28"; Value = "6" },
    @{ Cell = "G29"; Value = "6" },
    @{ Cell = "G30"; Value = "6" },
    @{ Cell = "G31"; Value = "6" },
    @{ Cell = "G32"; Value = "6" },
    @{ Cell = "G33"; Value = "6" },
    @{ Cell = "G34"; Value = "6" },
    @{ Cell = "G35"; Value = "6" },
    @{ Cell = "G36"; Value = "6" },
    @{ Cell = "G37"; Value = "6" },
    @{ Cell = "G38"; Value = "6" },
    @{ Cell = "D39"; Value = "0.02028" },
    @{ Cell = "E39"; Value = "6.48%" },
    @{ Cell = "G39"; Value = "6" },
    @{ Cell = "D40"; Value = "0.04977" },
    @{ Cell = "E40"; Value = "5.64%" },
    @{ Cell = "G40"; Value = "6" },
    @{ Cell = "D41"; Value = "0.007637" },
    @{ Cell = "E41"; Value = "0.57%" },
    @{ Cell = "G41"; Value = "6" },
    @{ Cell = "D42"; Value = "0.01027" },
    @{ Cell = "E42"; Value = "7.38%" },
    @{ Cell = "G42"; Value = "6" },
    @{ Cell = "D43"; Value = "0.1387" },
    @{ Cell = "E43"; Value = "4.59%" },
    @{ Cell = "G43"; Value = "6" },
    @{ Cell = "D44"; Value = "0.002103" },
    @{ Cell = "E44"; Value = "-0.42%" },
    @{ Cell = "G44"; Value = "6" },
    @{ Cell = "D45"; Value = "0.01196" },
    @{ Cell = "E45"; Value = "7.31%" },
    @{ Cell = "G45"; Value = "6" },
    @{ Cell = "D46"; Value = "0.00006442" },
    @{ Cell = "E46"; Value = "-0.53%" },
    @{ Cell = "G46"; Value = "6" },
    @{ Cell = "E47"; Value = "-0.05%" },
    @{ Cell = "G47"; Value = "6" },
    @{ Cell = "G48"; Value = "6" },
    @{ Cell = "E49"; Value = "-8.71%" },
    @{ Cell = "G49"; Value = "6" },
    @{ Cell = "D50"; Value = "0.00002100" },
    @{ Cell = "E50"; Value = "-0.05%" },
    @{ Cell = "G50"; Value = "6" },
    @{ Cell = "D51"; Value = "0.0002000" },
    @{ Cell = "E51"; Value = "-0.05%" },
    @{ Cell = "G51"; Value = "6" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
